$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (new) <= old Row 7 values
$ws.Range("D3").Value = 44229
$ws.Range("M3").Value = 55
$ws.Range("N3").Value = 11000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 11364
$ws.Range("S3").Value = 812

# Row 5 (new) <= old Row 3 values
$ws.Range("D5").Value = 44172
$ws.Range("M5").Value = 90
$ws.Range("N5").Value = 8500
$ws.Range("O5").Value = 9000
$ws.Range("P5").Value = 8806
$ws.Range("S5").Value = 629

# Row 6 (new) <= old Row 8 values
$ws.Range("D6").Value = 44181
$ws.Range("M6").Value = 65
$ws.Range("N6").Value = 9000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 9462
$ws.Range("S6").Value = 676

# Row 7 (new) <= old Row 6 values
$ws.Range("D7").Value = 44210
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 10000
$ws.Range("O7").Value = 11000
$ws.Range("P7").Value = 10357
$ws.Range("S7").Value = 740

# Row 8 (new) <= old Row 5 values
$ws.Range("D8").Value = 44232
$ws.Range("M8").Value = 60
$ws.Range("N8").Value = 11000
$ws.Range("O8").Value = 12000
$ws.Range("P8").Value = 11583
$ws.Range("S8").Value = 827
